$d = $word.ActiveDocument

function XmlEscape($s) {
    return $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

# Replaces the visible text of paragraph number $index (1-based, per
# $d.Paragraphs) via Range.InsertXML, leaving any leading empty run
# (e.g. stray <w:r/> artifacts) untouched and applying the given run
# properties (rPr) XML fragment to the new text run. Only applies the
# change if the paragraph's current text matches $oldText, as a guard
# against acting on the wrong paragraph.
function Set-ParagraphText($index, $oldText, $newText, $runPrXml) {
    $p = $d.Paragraphs.Item($index)
    $current = $p.Range.Text.TrimEnd("`r", "`a")
    if ($current -ne $oldText) {
        Write-Host "Set-ParagraphText: paragraph $index text mismatch; expected [$oldText] got [$current]"
        return
    }
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $escaped = XmlEscape($newText)
    $runXml = "<w:r>$runPrXml<w:t>$escaped</w:t></w:r>"
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# --- 1. Title (Heading1 "Play Lotus Luck Free: ..." and the identical bold
#          run repeated near the bottom). A plain whole-document Find/Replace
#          updates both occurrences in one shot without disturbing run
#          structure, since the Heading1 run has no leading empty run and the
#          bold run already differs in formatting from its own leading empty
#          run (so the engine keeps that boundary intact). ---
$d.Content.Find.Execute(
    "Play Lotus Luck Free: Detailed Review & Pros and Cons", $true, $true, `
    $false, $false, $false, $true, 1, $false, `
    "Play Lotus Luck Free - Review of Gameplay, Graphics, and More", 2) | Out-Null

# --- 2. "What we like" bullet list (paragraphs 37-40) ---
Set-ParagraphText 37 "Good compromise between win value and frequency of combinations" `
                      "1,024 ways to win" ""
Set-ParagraphText 38 "Detailed and well-designed graphics and symbols" `
                      "Good compromise between value of wins and frequency of winning combinations" ""
Set-ParagraphText 39 "1,024 ways to win with medium volatility and 95.62% RTP" `
                      "Detailed graphics and design" ""
Set-ParagraphText 40 "Relaxing gameplay with option for turbo mode and automatic spins" `
                      "Relaxing gameplay with no crashes or slowdowns" ""

# --- 3. "What we don't like" bullet list: update first bullet, delete the second ---
Set-ParagraphText 42 "Payment table lacks fluidity with no scrolling bar" `
                      "Lack of fluidity in payment table scrolling" ""

# Remove the "Limited number of themed symbols" paragraph entirely (para 43)
$p43 = $d.Paragraphs.Item(43)
if ($p43.Range.Text.TrimEnd("`r") -eq "Limited number of themed symbols") {
    $d.Range($p43.Range.Start, $p43.Range.End).Delete()
}

# --- 4. Italic description text near the bottom (now the last paragraph).
#          Scoping the Find to italic formatting makes the engine preserve
#          the leading empty run and the <w:i/> run properties. ---
$d.Content.Find.Font.Italic = 1
$d.Content.Find.Execute(
    "Discover the gameplay, features, and symbols of Lotus Luck online slot and find out what we like and don't like. Play for free and enjoy relaxing gameplay.", `
    $true, $true, $false, $false, $false, $true, 1, $false, `
    "Play Lotus Luck for free and enjoy its relaxing gameplay, detailed graphics, and 1,024 ways to win.", 2) | Out-Null
